# Insert a new weekly price record as row 237 on the single data sheet.
# This pushes every existing row from 237..336 down by one (to 238..337),
# which is exactly the shape of the target diff (each old D/J/K/L/M/.../P
# value reappears one row further down, and the sheet's used range grows
# from A1:R336 to A1:R337).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 237:336 down to 238:337, leaving a blank row at 237.
$ws.Rows("237:237").Insert()

# Populate the newly inserted row 237 with the new record.
$ws.Cells.Item(237, 1).Value  = 5
$ws.Cells.Item(237, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(237, 3).Value  = "Maule"
$ws.Cells.Item(237, 4).Value  = 44875
$ws.Cells.Item(237, 5).Value  = 7
$ws.Cells.Item(237, 6).Value  = 100112009
$ws.Cells.Item(237, 7).Value  = "Acelga"
$ws.Cells.Item(237, 8).Value  = "Sin especificar"
$ws.Cells.Item(237, 9).Value  = "Primera"
$ws.Cells.Item(237, 10).Value = 500
$ws.Cells.Item(237, 11).Value = 2000
$ws.Cells.Item(237, 12).Value = 2000
$ws.Cells.Item(237, 13).Value = 2000
$ws.Cells.Item(237, 14).Value = '$/docena de atados (4 kilos)'
$ws.Cells.Item(237, 15).Value = "Región del Maule"
$ws.Cells.Item(237, 16).Value = 500
$ws.Cells.Item(237, 17).Value = 4
$ws.Cells.Item(237, 18).Value = "Hortaliza"
